# filter finished / testing_copy_filter.py / v.8.8
# Append new log rows (160-164) to the bottom of the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records to append (Fecha, Hora, WC47 NACP, WC48 P5F, WC49 P5H, WV50 FILTER, SPL)
$newRows = @(
    @("2024-05-14", "12:30:16", "Fallo fijador tapa", "-", "-", "-", "-"),
    @("2024-05-14", "12:30:23", "-", "No detecta presencia power CP", "-", "-", "-"),
    @("2024-05-14", "12:34:20", "Ascensor no sube", "-", "-", "-", "-"),
    @("2024-05-14", "12:34:24", "Etiquetadora", "-", "-", "-", "-"),
    @("2024-05-14", "12:34:28", "Fallo fijador tapa", "-", "-", "-", "-")
)

$startRow = 160

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Fecha column looks like a date, so prefix with an apostrophe to force
    # Excel to store it as literal text (matches the other rows in the sheet).
    $ws.Range("A$r").Value = "'" + $rowData[0]
    $ws.Range("B$r").Value = $rowData[1]
    $ws.Range("C$r").Value = $rowData[2]
    $ws.Range("D$r").Value = $rowData[3]
    $ws.Range("E$r").Value = $rowData[4]
    $ws.Range("F$r").Value = $rowData[5]
    $ws.Range("G$r").Value = $rowData[6]
}
